$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in column H, matching the formatting of the
# existing header cells (e.g. G1 -> bold, bordered, centered style).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# Add the corresponding data value for the new column
$ws.Range("H2").Value = 0
